$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 1.440985666666667
$ws.Range("H2").Value = 4.322957000000001
$ws.Range("I2").Value = 0.1098365531732288
$ws.Range("J2").Value = 0.1230162332390494
$ws.Range("M2").Value = 13.89934866666667
$ws.Range("N2").Value = 41.69804600000001
$ws.Range("O2").Value = 0.04853507553134179
$ws.Range("P2").Value = 0.04999273878390351
$ws.Range("Q2").Value = 20.02876220466912
$ws.Range("R2").Value = 180.258859842022
$ws.Range("S2").Value = 0.005330925404364898
$ws.Range("T2").Value = 0.006149918414499543

# Row 3
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 1.440985666666667
$ws.Range("H3").Value = 4.322957000000001
$ws.Range("I3").Value = 0.1098365531732288
$ws.Range("J3").Value = 0.1230162332390494
$ws.Range("O3").Value = 0.245697991654417
$ws.Range("P3").Value = 0.253077086664408
$ws.Range("Q3").Value = 101.3911402246261
$ws.Range("R3").Value = 912.5202620216351
$ws.Range("S3").Value = 0.0269866205249059
$ws.Range("T3").Value = 0.03113258992056792

# Row 4
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 1.440985666666667
$ws.Range("H4").Value = 4.322957000000001
$ws.Range("I4").Value = 0.1098365531732288
$ws.Range("J4").Value = 0.1230162332390494
$ws.Range("M4").Value = 82.007665
$ws.Range("N4").Value = 246.022995
$ws.Range("O4").Value = 0.2863622109480123
$ws.Range("P4").Value = 0.2949625822722868
$ws.Range("Q4").Value = 118.1718698218017
$ws.Range("R4").Value = 1063.546828396215
$ws.Range("S4").Value = 0.03145303820959472
$ws.Range("T4").Value = 0.03628518581759992

# Row 5
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 1.440985666666667
$ws.Range("H5").Value = 4.322957000000001
$ws.Range("I5").Value = 0.1098365531732288
$ws.Range("J5").Value = 0.1230162332390494
$ws.Range("M5").Value = 25.0501465
$ws.Range("N5").Value = 50.100293
$ws.Range("O5").Value = 0.0874724982879541
$ws.Range("P5").Value = 0.06006638442832619
$ws.Range("Q5").Value = 36.09690205440018
$ws.Range("R5").Value = 216.581412326401
$ws.Range("S5").Value = 0.009607677709400035
$ws.Range("T5").Value = 0.007389140356661377

# Row 6
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 1.440985666666667
$ws.Range("H6").Value = 4.322957000000001
$ws.Range("I6").Value = 0.1098365531732288
$ws.Range("J6").Value = 0.1230162332390494
$ws.Range("M6").Value = 95.05788666666668
$ws.Range("N6").Value = 285.17366
$ws.Range("O6").Value = 0.3319322235782747
$ws.Range("P6").Value = 0.3419012078510756
$ws.Range("Q6").Value = 136.9770521902911
$ws.Range("R6").Value = 1232.79346971262
$ws.Range("S6").Value = 0.03645829132496323
$ws.Range("T6").Value = 0.04205939872972062

# Row 7
$ws.Range("I7").Value = 0.5687502547919595
$ws.Range("J7").Value = 0.6369966279614609
$ws.Range("M7").Value = 13.89934866666667
$ws.Range("N7").Value = 41.69804600000001
$ws.Range("O7").Value = 0.04853507553134179
$ws.Range("P7").Value = 0.04999273878390351
$ws.Range("Q7").Value = 103.7119545176116
$ws.Range("R7").Value = 933.407590658504
$ws.Range("S7").Value = 0.02760433657479764
$ws.Range("T7").Value = 0.03184520602790468

# Row 8
$ws.Range("I8").Value = 0.5687502547919595
$ws.Range("J8").Value = 0.6369966279614609
$ws.Range("O8").Value = 0.245697991654417
$ws.Range("P8").Value = 0.253077086664408
$ws.Range("S8").Value = 0.1397407953553224
$ws.Range("T8").Value = 0.1612092508195383

# Row 9
$ws.Range("I9").Value = 0.5687502547919595
$ws.Range("J9").Value = 0.6369966279614609
$ws.Range("M9").Value = 82.007665
$ws.Range("N9").Value = 246.022995
$ws.Range("O9").Value = 0.2863622109480123
$ws.Range("P9").Value = 0.2949625822722868
$ws.Range("Q9").Value = 611.9117828141533
$ws.Range("R9").Value = 5507.20604532738
$ws.Range("S9").Value = 0.1628685804394709
$ws.Range("T9").Value = 0.1878901702822517

# Row 10
$ws.Range("I10").Value = 0.5687502547919595
$ws.Range("J10").Value = 0.6369966279614609
$ws.Range("M10").Value = 25.0501465
$ws.Range("N10").Value = 50.100293
$ws.Range("O10").Value = 0.0874724982879541
$ws.Range("P10").Value = 0.06006638442832619
$ws.Range("Q10").Value = 186.9152085304553
$ws.Range("R10").Value = 1121.491251182732
$ws.Range("S10").Value = 0.04975000568856314
$ws.Range("T10").Value = 0.03826208433468059

# Row 11
$ws.Range("I11").Value = 0.5687502547919595
$ws.Range("J11").Value = 0.6369966279614609
$ws.Range("M11").Value = 95.05788666666668
$ws.Range("N11").Value = 285.17366
$ws.Range("O11").Value = 0.3319322235782747
$ws.Range("P11").Value = 0.3419012078510756
$ws.Range("Q11").Value = 709.2878562113156
$ws.Range("R11").Value = 6383.59070590184
$ws.Range("S11").Value = 0.1887865367338054
$ws.Range("T11").Value = 0.2177899164970858

# Row 12
$ws.Range("G12").Value = 4.2167365
$ws.Range("H12").Value = 8.433472999999999
$ws.Range("I12").Value = 0.3214131920348118
$ws.Range("J12").Value = 0.2399871387994896
$ws.Range("M12").Value = 13.89934866666667
$ws.Range("N12").Value = 41.69804600000001
$ws.Range("O12").Value = 0.04853507553134179
$ws.Range("P12").Value = 0.04999273878390351
$ws.Range("Q12").Value = 58.60989084895967
$ws.Range("R12").Value = 351.659345093758
$ws.Range("S12").Value = 0.01559981355217925
$ws.Range("T12").Value = 0.01199761434149928

# Row 13
$ws.Range("G13").Value = 4.2167365
$ws.Range("H13").Value = 8.433472999999999
$ws.Range("I13").Value = 0.3214131920348118
$ws.Range("J13").Value = 0.2399871387994896
$ws.Range("O13").Value = 0.245697991654417
$ws.Range("P13").Value = 0.253077086664408
$ws.Range("Q13").Value = 296.6994964986691
$ws.Range("R13").Value = 1780.196978992015
$ws.Range("S13").Value = 0.07897057577418871
$ws.Range("T13").Value = 0.06073524592430174

# Row 14
$ws.Range("G14").Value = 4.2167365
$ws.Range("H14").Value = 8.433472999999999
$ws.Range("I14").Value = 0.3214131920348118
$ws.Range("J14").Value = 0.2399871387994896
$ws.Range("M14").Value = 82.007665
$ws.Range("N14").Value = 246.022995
$ws.Range("O14").Value = 0.2863622109480123
$ws.Range("P14").Value = 0.2949625822722868
$ws.Range("Q14").Value = 345.8047142852725
$ws.Range("R14").Value = 2074.828285711635
$ws.Range("S14").Value = 0.09204059229894677
$ws.Range("T14").Value = 0.07078722617243516

# Row 15
$ws.Range("G15").Value = 4.2167365
$ws.Range("H15").Value = 8.433472999999999
$ws.Range("I15").Value = 0.3214131920348118
$ws.Range("J15").Value = 0.2399871387994896
$ws.Range("M15").Value = 25.0501465
$ws.Range("N15").Value = 50.100293
$ws.Range("O15").Value = 0.0874724982879541
$ws.Range("P15").Value = 0.06006638442832619
$ws.Range("Q15").Value = 105.6298670768972
$ws.Range("R15").Value = 422.519468307589
$ws.Range("S15").Value = 0.02811481488999094
$ws.Range("T15").Value = 0.01441515973698422

# Row 16
$ws.Range("G16").Value = 4.2167365
$ws.Range("H16").Value = 8.433472999999999
$ws.Range("I16").Value = 0.3214131920348118
$ws.Range("J16").Value = 0.2399871387994896
$ws.Range("M16").Value = 95.05788666666668
$ws.Range("N16").Value = 285.17366
$ws.Range("O16").Value = 0.3319322235782747
$ws.Range("P16").Value = 0.3419012078510756
$ws.Range("Q16").Value = 400.8340603201967
$ws.Range("R16").Value = 2405.00436192118
$ws.Range("S16").Value = 0.1066873955195061
$ws.Range("T16").Value = 0.08205189262426923
